# Auto-generated edit script: updates date formatting (YYYY-MM-DD -> YYYY.MM.DD)
# and refreshes "想去人数" (attendance) counts across the "展览", "演出" and
# "全部类型" worksheets, matching the upstream gh-pages data refresh at 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws_sheet1 = $wb.Worksheets.Item("展览")

# Force text format on the date column so "YYYY.MM.DD" strings are not
# auto-converted into date serial numbers by the COM layer.
$ws_sheet1.Range("B2:B28").NumberFormat = "@"

# Update start-date text values
$ws_sheet1.Range("B2").Value = "2024.02.08"
$ws_sheet1.Range("B3").Value = "2024.02.14"
$ws_sheet1.Range("B4").Value = "2024.02.14"
$ws_sheet1.Range("B5").Value = "2024.02.14"
$ws_sheet1.Range("B6").Value = "2024.02.16"
$ws_sheet1.Range("B7").Value = "2024.02.16"
$ws_sheet1.Range("B8").Value = "2024.02.16"
$ws_sheet1.Range("B9").Value = "2024.02.24"
$ws_sheet1.Range("B10").Value = "2024.02.24"
$ws_sheet1.Range("B11").Value = "2024.02.24"
$ws_sheet1.Range("B12").Value = "2024.02.25"
$ws_sheet1.Range("B13").Value = "2024.02.25"
$ws_sheet1.Range("B14").Value = "2024.03.08"
$ws_sheet1.Range("B15").Value = "2024.03.17"
$ws_sheet1.Range("B16").Value = "2024.03.23"
$ws_sheet1.Range("B17").Value = "2024.04.04"
$ws_sheet1.Range("B18").Value = "2024.04.06"
$ws_sheet1.Range("B19").Value = "2024.04.13"
$ws_sheet1.Range("B20").Value = "2024.04.21"
$ws_sheet1.Range("B21").Value = "2024.05.01"
$ws_sheet1.Range("B22").Value = "2024.05.01"
$ws_sheet1.Range("B23").Value = "2024.05.02"
$ws_sheet1.Range("B24").Value = "2024.05.02"
$ws_sheet1.Range("B25").Value = "2024.05.02"
$ws_sheet1.Range("B26").Value = "2024.05.02"
$ws_sheet1.Range("B27").Value = "2024.05.03"
$ws_sheet1.Range("B28").Value = "2024.05.03"

# Update "想去人数" (want-to-go count) values
$ws_sheet1.Range("F3").Value = 1438
$ws_sheet1.Range("F7").Value = 12018
$ws_sheet1.Range("F8").Value = 4459
$ws_sheet1.Range("F9").Value = 38
$ws_sheet1.Range("F10").Value = 58
$ws_sheet1.Range("F13").Value = 2579
$ws_sheet1.Range("F14").Value = 1115
$ws_sheet1.Range("F15").Value = 176
$ws_sheet1.Range("F16").Value = 58
$ws_sheet1.Range("F17").Value = 5202
$ws_sheet1.Range("F21").Value = 11406
$ws_sheet1.Range("F22").Value = 11429
$ws_sheet1.Range("F24").Value = 54

# --- Sheet "演出" ---
$ws_sheet2 = $wb.Worksheets.Item("演出")

# Force text format on the date column so "YYYY.MM.DD" strings are not
# auto-converted into date serial numbers by the COM layer.
$ws_sheet2.Range("B2:B2").NumberFormat = "@"

# Update start-date text values
$ws_sheet2.Range("B2").Value = "2024.03.03"

# --- Sheet "全部类型" ---
$ws_sheet4 = $wb.Worksheets.Item("全部类型")

# Force text format on the date column so "YYYY.MM.DD" strings are not
# auto-converted into date serial numbers by the COM layer.
$ws_sheet4.Range("B2:B29").NumberFormat = "@"

# Update start-date text values
$ws_sheet4.Range("B2").Value = "2024.02.08"
$ws_sheet4.Range("B3").Value = "2024.02.14"
$ws_sheet4.Range("B4").Value = "2024.02.14"
$ws_sheet4.Range("B5").Value = "2024.02.14"
$ws_sheet4.Range("B6").Value = "2024.02.16"
$ws_sheet4.Range("B7").Value = "2024.02.16"
$ws_sheet4.Range("B8").Value = "2024.02.16"
$ws_sheet4.Range("B9").Value = "2024.02.24"
$ws_sheet4.Range("B10").Value = "2024.02.24"
$ws_sheet4.Range("B11").Value = "2024.02.24"
$ws_sheet4.Range("B12").Value = "2024.02.25"
$ws_sheet4.Range("B13").Value = "2024.02.25"
$ws_sheet4.Range("B14").Value = "2024.03.03"
$ws_sheet4.Range("B15").Value = "2024.03.08"
$ws_sheet4.Range("B16").Value = "2024.03.17"
$ws_sheet4.Range("B17").Value = "2024.03.23"
$ws_sheet4.Range("B18").Value = "2024.04.04"
$ws_sheet4.Range("B19").Value = "2024.04.06"
$ws_sheet4.Range("B20").Value = "2024.04.13"
$ws_sheet4.Range("B21").Value = "2024.04.21"
$ws_sheet4.Range("B22").Value = "2024.05.01"
$ws_sheet4.Range("B23").Value = "2024.05.01"
$ws_sheet4.Range("B24").Value = "2024.05.02"
$ws_sheet4.Range("B25").Value = "2024.05.02"
$ws_sheet4.Range("B26").Value = "2024.05.02"
$ws_sheet4.Range("B27").Value = "2024.05.02"
$ws_sheet4.Range("B28").Value = "2024.05.03"
$ws_sheet4.Range("B29").Value = "2024.05.03"

# Update "想去人数" (want-to-go count) values
$ws_sheet4.Range("F3").Value = 1438
$ws_sheet4.Range("F7").Value = 12018
$ws_sheet4.Range("F8").Value = 4459
$ws_sheet4.Range("F9").Value = 38
$ws_sheet4.Range("F10").Value = 58
$ws_sheet4.Range("F13").Value = 2579
$ws_sheet4.Range("F15").Value = 1115
$ws_sheet4.Range("F16").Value = 176
$ws_sheet4.Range("F17").Value = 58
$ws_sheet4.Range("F18").Value = 5202
$ws_sheet4.Range("F22").Value = 11406
$ws_sheet4.Range("F23").Value = 11429
$ws_sheet4.Range("F25").Value = 54

